$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting existing rows 108..228 down to 109..229
$ws.Rows(108).Insert()

# Populate the newly inserted row 108 with the new record's data
$ws.Cells.Item(108, 1).Value = 10
$ws.Cells.Item(108, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(108, 3).Value = "La Araucanía"
$ws.Cells.Item(108, 4).Value = 44789
$ws.Cells.Item(108, 5).Value = 9
$ws.Cells.Item(108, 6).Value = 100112005
$ws.Cells.Item(108, 7).Value = "Puerro"
$ws.Cells.Item(108, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 65
$ws.Cells.Item(108, 11).Value = 16000
$ws.Cells.Item(108, 12).Value = 16000
$ws.Cells.Item(108, 13).Value = 16000
$ws.Cells.Item(108, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(108, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(108, 16).Value = 1333
$ws.Cells.Item(108, 17).Value = 12
$ws.Cells.Item(108, 18).Value = "Hortaliza"
